$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.565.01"
Set-TextValue "E2" "  +3.74%  "
Set-TextValue "D3" "1.848.21"
Set-TextValue "E3" "  +2.91%  "
Set-TextValue "D4" "1.030"
Set-TextValue "E4" "  +2.91%  "
Set-TextValue "D5" "319.45"
Set-TextValue "E5" "  +3.33%  "
Set-TextValue "D6" "1.025"
Set-TextValue "E6" "  +2.50%  "
Set-TextValue "D7" "0.4368"
Set-TextValue "E7" "  +1.99%  "
Set-TextValue "D8" "0.3749"
Set-TextValue "E8" "  +3.36%  "
Set-TextValue "D9" "0.07400"
Set-TextValue "E9" "  +3.10%  "
Set-TextValue "D10" "0.8772"
Set-TextValue "E10" "  +2.32%  "
Set-TextValue "D11" "21.53"
Set-TextValue "E11" "  +4.32%  "
Set-TextValue "D12" "1.860.39"
Set-TextValue "E12" "  -1.39%  "
Set-TextValue "D13" "5.486"
Set-TextValue "E13" "  +3.68%  "
Set-TextValue "D14" "6.685"
Set-TextValue "E14" "  +2.34%  "
Set-TextValue "E15" "  +3.53%  "
Set-TextValue "D16" "82.78"
Set-TextValue "E16" "  +3.77%  "
Set-TextValue "D17" "1.032"
Set-TextValue "E17" "  +3.13%  "
Set-TextValue "D18" "0.000009016"
Set-TextValue "E18" "  +2.72%  "
Set-TextValue "D19" "1.025"
Set-TextValue "E19" "  +2.27%  "
Set-TextValue "D20" "15.43"
Set-TextValue "E20" "  +2.55%  "
Set-TextValue "D21" "27.583.09"
Set-TextValue "E21" "  +3.78%  "
Set-TextValue "D22" "5.253"
Set-TextValue "E22" "  +2.20%  "
Set-TextValue "D23" "11.21"
Set-TextValue "E23" "  +1.53%  "
Set-TextValue "D24" "2.071.47"
Set-TextValue "E24" "  -1.52%  "
Set-TextValue "D25" "157.33"
Set-TextValue "E25" "  +3.56%  "
Set-TextValue "D26" "1.926"
Set-TextValue "E26" "  +5.62%  "
Set-TextValue "D27" "18.73"
Set-TextValue "E27" "  +3.05%  "
Set-TextValue "D28" "5.253"
Set-TextValue "E28" "  +1.68%  "
Set-TextValue "D29" "1.949"
Set-TextValue "E29" "  +3.24%  "
Set-TextValue "D30" "116.04"
Set-TextValue "D31" "0.09078"
Set-TextValue "E31" "  +1.97%  "
Set-TextValue "D32" "1.207"
Set-TextValue "E32" "  +5.06%  "
Set-TextValue "E33" "  +2.65%  "
Set-TextValue "D34" "4.502"
Set-TextValue "E34" "  +3.10%  "
Set-TextValue "D35" "2.870"
Set-TextValue "E35" "  +4.72%  "
Set-TextValue "D36" "1.027"
Set-TextValue "E36" "  +2.42%  "
Set-TextValue "D37" "1.144"
Set-TextValue "E37" "  +3.06%  "
Set-TextValue "D38" "0.01977"
Set-TextValue "E38" "  +3.97%  "
Set-TextValue "D39" "0.05268"
Set-TextValue "E39" "  +2.11%  "
Set-TextValue "D40" "0.5171"
Set-TextValue "E40" "  +3.57%  "
Set-TextValue "D41" "2.794"
Set-TextValue "E41" "  +7.19%  "
Set-TextValue "D42" "0.1674"
Set-TextValue "E42" "  +3.09%  "
Set-TextValue "D43" "6.699"
Set-TextValue "E43" "  +3.97%  "
Set-TextValue "D44" "8.546"
Set-TextValue "E44" "  +3.68%  "
Set-TextValue "B45" "Quant"
Set-TextValue "C45" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D45" "108.83"
Set-TextValue "E45" "  +2.76%  "
Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "10.56"
Set-TextValue "E46" "  +2.46%  "
Set-TextValue "E47" "  +4.30%  "
Set-TextValue "D48" "0.4652"
Set-TextValue "E48" "  +3.15%  "
Set-TextValue "D49" "0.06368"
Set-TextValue "E49" "  +2.53%  "
Set-TextValue "D50" "1.894"
Set-TextValue "E50" "  +6.30%  "
Set-TextValue "D51" "39.53"
Set-TextValue "E51" "  +6.62%  "
